# Update absenteeism data rows 2-11 with new values per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(49638, "Sra. Júlia da Rocha",        "Vendas",            "Viagem de negócios", 7, 45091, 11476.71),
    @(41594, "Srta. Heloísa Nascimento",   "Recursos Humanos",  "Outros",              5, 45086, 10107.15),
    @(73189, "Stephany Duarte",             "Jurídico",          "Consulta médica",    1, 45105, 6941.61),
    @(879,   "Fernando Nogueira",           "Jurídico",          "Viagem de negócios", 6, 45092, 5893.07),
    @(41172, "Helena Pires",                "Marketing",         "Consulta médica",    6, 45101, 4505.12),
    @(74310, "Stella Cardoso",              "P&D",               "Viagem de negócios", 3, 45093, 9814.94),
    @(735,   "Guilherme Jesus",             "Vendas",            "Viagem de negócios", 2, 45101, 7620.48),
    @(51789, "Júlia Ribeiro",               "Marketing",         "Consulta médica",    7, 45083, 10621.78),
    @(19705, "Ana Clara da Rosa",           "Engenharia",        "Doença",             4, 45080, 5546.39),
    @(35991, "Sra. Mariane Freitas",        "Vendas",            "Outros",             5, 45104, 3070.78)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]

    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
}
